$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# ---- New text values (replace the two handed-off source files with the
#      new consolidated handoff file + renamed dependency file, and flip the
#      status from "Handed back: in sync with en-US" to "Ready for handoff") ----
$newMd1    = "cfd858bf-bc47-4158-83a6-a0a57bb5e0d2.md"
$newMd2    = "ffffe695b2f2-bd49-4971-9823-e3e11da8248d.md"
$newStatus = "Ready for handoff"

$zhXlf  = "cfd858bf-bc47-4158-83a6-a0a57bb5e0d2.da0aa32d138b3d26f96f0e8df8019a170cfd18d5.zh-cn.xlf"
$deXlf  = "cfd858bf-bc47-4158-83a6-a0a57bb5e0d2.da0aa32d138b3d26f96f0e8df8019a170cfd18d5.de-de.xlf"
$zhDate = "2016-03-10 05:22:17"
$deDate = "2016-03-10 05:22:22"
$epoch  = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Overview sheet: only the cell text changes (same shared strings get their
# text swapped in place because every other usage below is updated too).
# ---------------------------------------------------------------------------
$wsOverview.Range("A2").Value = $newMd1
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus

$wsOverview.Range("A3").Value = $newMd2
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh.Range("A2").Value = $newMd1
$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("C2").Value = $zhXlf
$wsZh.Range("D2").Value = $zhDate
$wsZh.Range("G2").Value = $epoch
$wsZh.Range("H2").Value = "Include"

$wsZh.Range("A3").Value = $newMd2
$wsZh.Range("B3").Value = $newStatus
$wsZh.Range("C3").Value = $zhXlf
$wsZh.Range("D3").Value = $zhDate
$wsZh.Range("G3").Value = $epoch
$wsZh.Range("H3").Value = "Include"

# drop the now-unused Latest Target File / Latest Handback File columns
$wsZh.Range("E2:F3").Clear()

# row 4 (.localization-config) keeps its text but the shared-string index
# shifts once the orphaned strings above are garbage collected.
$wsZh.Range("D4").Value = $epoch
$wsZh.Range("G4").Value = $epoch
$wsZh.Range("H4").Value = "Ignored"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe.Range("A2").Value = $newMd1
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("C2").Value = $deXlf
$wsDe.Range("D2").Value = $deDate
$wsDe.Range("G2").Value = $epoch
$wsDe.Range("H2").Value = "Include"

$wsDe.Range("A3").Value = $newMd2
$wsDe.Range("B3").Value = $newStatus
$wsDe.Range("C3").Value = $deXlf
$wsDe.Range("D3").Value = $deDate
$wsDe.Range("G3").Value = $epoch
$wsDe.Range("H3").Value = "Include"

$wsDe.Range("E2:F3").Clear()

$wsDe.Range("D4").Value = $epoch
$wsDe.Range("G4").Value = $epoch
$wsDe.Range("H4").Value = "Ignored"

# ---------------------------------------------------------------------------
# Hyperlinks: rebuild zh-cn / de-de link tables without the removed E/F
# columns (Overview's links are untouched, matching the unchanged diff there).
# ---------------------------------------------------------------------------
$wsZh.Range("A1").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/beacc8af0d380dda163a7cc9282aaedaac687098/e2e/$newMd1", "", "", $newMd1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/02417d0f3f3d8a047caa8a6e40c4f83d533ac06d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf", "", "", $zhXlf) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/beacc8af0d380dda163a7cc9282aaedaac687098/e2e/$newMd2", "", "", $newMd2) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/02417d0f3f3d8a047caa8a6e40c4f83d533ac06d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlf", "", "", $zhXlf) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/beacc8af0d380dda163a7cc9282aaedaac687098/.localization-config", "", "", ".localization-config") | Out-Null

$wsDe.Range("A1").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/beacc8af0d380dda163a7cc9282aaedaac687098/e2e/$newMd1", "", "", $newMd1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a0d86c80d6da183b731ca0fb9147aa182189a663/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf", "", "", $deXlf) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/beacc8af0d380dda163a7cc9282aaedaac687098/e2e/$newMd2", "", "", $newMd2) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a0d86c80d6da183b731ca0fb9147aa182189a663/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlf", "", "", $deXlf) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/beacc8af0d380dda163a7cc9282aaedaac687098/.localization-config", "", "", ".localization-config") | Out-Null
